$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "audioC" header text to "audioX"
$ws.Range("C1").Value2 = "audioX"

# Swap columns B and C in the header row so the renamed "audioX" column
# moves to B1 and "audioB" moves to C1 (fixes naming/order of post test column)
$colB = $ws.Range("B1").Value2
$colC = $ws.Range("C1").Value2
$ws.Range("B1").Value2 = $colC
$ws.Range("C1").Value2 = $colB

# Update the active selection to match the saved view state
$ws.Range("M32").Select()
